# Update "想去人数" (interest count) figures for two events that appear
# on both the "展览" (Exhibition) sheet and the "全部类型" (All Types) sheet.
#
#   - 南宁·草莓动漫节               : 1213 -> 1218
#   - 南宁·第一届ANE·DACG动漫嘉年华 : 608  -> 610

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 4 is 草莓动漫节 (F4), row 5 is ANE·DACG动漫嘉年华 (F5)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 1218
$wsExhibition.Range("F5").Value = 610

# Sheet "全部类型": row 4 is 草莓动漫节 (F4), row 6 is ANE·DACG动漫嘉年华 (F6)
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F4").Value = 1218
$wsAllTypes.Range("F6").Value = 610
